# Added more test data.
# Sheet "Test": fix up row 8 (new VIN / mileage / used-vehicle program) and
# append rows 9-13 with the same kind of automation test data.
# Sheet "Two_program": append rows 7-9, mirroring the Test-sheet rows that
# have a second program column (F).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Test")
$ws2 = $wb.Worksheets.Item("Two_program")

$newVin = "5J6RW2H89NA004619"

# ---- Test!A8:J8 - correct the existing last row ----------------------
$ws1.Range("C8").Value = $newVin
$ws1.Range("D8").Value = 2345
$ws1.Range("E8").Value = "Used Vehicle - RNL"

# Give row 8 the "pasted from web" look (MuseoSans-300) that the rest of
# the newly entered rows pick up too.
$ws1.Range("A8:J8").Font.Name = "MuseoSans-300"
$ws1.Range("A8:J8").Font.Size = 10
$ws1.Range("A8:J8").Font.Color = 3355443

# ---- Two_program!A7:K7 - new row paired with Test!A8:J8 ---------------
$ws2.Range("A7").Value = "Automation 8"
$ws2.Range("B7").Value = "Test 8"
$ws2.Range("C7").Value = $newVin
$ws2.Range("D7").Value = 2345
$ws2.Range("E7").Value = "Used Vehicle - RNL"
$ws2.Range("F7").Value = "Used Vehicle - SNI"
$ws2.Range("G7").Value = "one"
$ws2.Range("H7").Value = "one"
$ws2.Range("I7").Value = "Auto Test"
$ws2.Range("J7").Value = 45678
$ws2.Range("K7").Value = 9999999

$ws2.Range("C7").Font.Name = "MuseoSans-300"
$ws2.Range("C7").Font.Size = 10
$ws2.Range("C7").Font.Color = 3355443

# ---- Test!A9:J9 --------------------------------------------------------
$ws1.Range("A9").Value = "Automation 9"
$ws1.Range("B9").Value = "Test 9"
$ws1.Range("C9").Value = $newVin
$ws1.Range("D9").Value = 2345
$ws1.Range("E9").Value = "Used Vehicle - SNE"
$ws1.Range("F9").Value = "one"
$ws1.Range("G9").Value = "one"
$ws1.Range("H9").Value = "Auto Test"
$ws1.Range("I9").Value = 45678
$ws1.Range("J9").Value = 9999999

$ws1.Range("A9:J9").Font.Name = "MuseoSans-300"
$ws1.Range("A9:J9").Font.Size = 10
$ws1.Range("A9:J9").Font.Color = 3355443

# ---- Two_program!A8:K8 --------------------------------------------------
$ws2.Range("A8").Value = "Automation 9"
$ws2.Range("B8").Value = "Test 9"
$ws2.Range("C8").Value = $newVin
$ws2.Range("D8").Value = 2345
$ws2.Range("E8").Value = "Used Vehicle - SNE"
$ws2.Range("F8").Value = "Used Vehicle - SNL"
$ws2.Range("G8").Value = "one"
$ws2.Range("H8").Value = "one"
$ws2.Range("I8").Value = "Auto Test"
$ws2.Range("J8").Value = 45678
$ws2.Range("K8").Value = 9999999

# ---- Test!A10:J10 -------------------------------------------------------
$ws1.Range("A10").Value = "Automation 10"
$ws1.Range("B10").Value = "Test 10"
$ws1.Range("C10").Value = $newVin
$ws1.Range("D10").Value = 2345
$ws1.Range("E10").Value = "Used Vehicle - SNF"
$ws1.Range("F10").Value = "one"
$ws1.Range("G10").Value = "one"
$ws1.Range("H10").Value = "Auto Test"
$ws1.Range("I10").Value = 45678
$ws1.Range("J10").Value = 9999999

$ws1.Range("A10:J10").Font.Name = "MuseoSans-300"
$ws1.Range("A10:J10").Font.Size = 10
$ws1.Range("A10:J10").Font.Color = 3355443

# ---- Two_program!A9:K9 --------------------------------------------------
$ws2.Range("A9").Value = "Automation 10"
$ws2.Range("B9").Value = "Test 10"
$ws2.Range("C9").Value = $newVin
$ws2.Range("D9").Value = 2345
$ws2.Range("E9").Value = "Used Vehicle - SNF"
$ws2.Range("F9").Value = "Used Vehicle - SPK"
$ws2.Range("G9").Value = "one"
$ws2.Range("H9").Value = "one"
$ws2.Range("I9").Value = "Auto Test"
$ws2.Range("J9").Value = 45678
$ws2.Range("K9").Value = 9999999

# ---- Test!A11:J11 -------------------------------------------------------
$ws1.Range("A11").Value = "Automation 11"
$ws1.Range("B11").Value = "Test 11"
$ws1.Range("C11").Value = $newVin
$ws1.Range("D11").Value = 2345
$ws1.Range("E11").Value = "Used Vehicle - SNI"
$ws1.Range("F11").Value = "one"
$ws1.Range("G11").Value = "one"
$ws1.Range("H11").Value = "Auto Test"
$ws1.Range("I11").Value = 45678
$ws1.Range("J11").Value = 9999999

$ws1.Range("A11:J11").Font.Name = "MuseoSans-300"
$ws1.Range("A11:J11").Font.Size = 10
$ws1.Range("A11:J11").Font.Color = 3355443

# ---- Test!A12:J12 -------------------------------------------------------
$ws1.Range("A12").Value = "Automation 12"
$ws1.Range("B12").Value = "Test 12"
$ws1.Range("C12").Value = $newVin
$ws1.Range("D12").Value = 2345
$ws1.Range("E12").Value = "Used Vehicle - SNL"
$ws1.Range("F12").Value = "one"
$ws1.Range("G12").Value = "one"
$ws1.Range("H12").Value = "Auto Test"
$ws1.Range("I12").Value = 45678
$ws1.Range("J12").Value = 9999999

$ws1.Range("A12:J12").Font.Name = "MuseoSans-300"
$ws1.Range("A12:J12").Font.Size = 10
$ws1.Range("A12:J12").Font.Color = 3355443

# ---- Test!A13:J13 -------------------------------------------------------
$ws1.Range("A13").Value = "Automation 13"
$ws1.Range("B13").Value = "Test 13"
$ws1.Range("C13").Value = $newVin
$ws1.Range("D13").Value = 2345
$ws1.Range("E13").Value = "Used Vehicle - SPK"
$ws1.Range("F13").Value = "one"
$ws1.Range("G13").Value = "one"
$ws1.Range("H13").Value = "Auto Test"
$ws1.Range("I13").Value = 45678
$ws1.Range("J13").Value = 9999999

$ws1.Range("A13:J13").Font.Name = "MuseoSans-300"
$ws1.Range("A13:J13").Font.Size = 10
$ws1.Range("A13:J13").Font.Color = 3355443

# ---- selection / scroll position, matching where the user ended up -----
$ws2.Activate()
$ws2.Range("A10").Select() | Out-Null

$ws1.Activate()
$ws1.Range("B1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("D14").Select() | Out-Null

Write-Host "Added more test data"
